$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "Questions"

# Header row: swap IS STATISTICAL / IS HISTORIC columns (F/G)
$ws.Cells.Item(1,6).Value = 'IS HISTORIC'
$ws.Cells.Item(1,7).Value = 'IS STATISTICAL'

# New column I: difficulty rating for each question
$ws.Cells.Item(2,9).Value = 'easy'
$ws.Cells.Item(3,9).Value = 'easy'
$ws.Cells.Item(4,9).Value = 'difficult'
$ws.Cells.Item(5,9).Value = 'difficult'
$ws.Cells.Item(6,9).Value = 'difficult'
$ws.Cells.Item(7,9).Value = 'easy'
$ws.Cells.Item(8,9).Value = 'easy'
$ws.Cells.Item(9,9).Value = 'easy'
$ws.Cells.Item(10,9).Value = 'difficult'
$ws.Cells.Item(11,9).Value = 'extreme'
$ws.Cells.Item(12,9).Value = 'extreme'
$ws.Cells.Item(13,9).Value = 'extreme'
$ws.Cells.Item(14,9).Value = 'easy'
$ws.Cells.Item(15,9).Value = 'difficult'
$ws.Cells.Item(16,9).Value = 'extreme'
$ws.Cells.Item(17,9).Value = 'easy'
$ws.Cells.Item(18,9).Value = 'difficult'
$ws.Cells.Item(19,9).Value = 'difficult'
$ws.Cells.Item(20,9).Value = 'easy'
$ws.Cells.Item(21,9).Value = 'difficult'
$ws.Cells.Item(22,9).Value = 'easy'
$ws.Cells.Item(23,9).Value = 'easy'
$ws.Cells.Item(24,9).Value = 'difficult'
$ws.Cells.Item(25,9).Value = 'easy'
$ws.Cells.Item(26,9).Value = 'difficult'
$ws.Cells.Item(27,9).Value = 'difficult'
$ws.Cells.Item(28,9).Value = 'easy'
$ws.Cells.Item(29,9).Value = 'difficult'
$ws.Cells.Item(30,9).Value = 'easy'
$ws.Cells.Item(31,9).Value = 'easy'
$ws.Cells.Item(32,9).Value = 'easy'
$ws.Cells.Item(33,9).Value = 'easy'
$ws.Cells.Item(34,9).Value = 'easy'
$ws.Cells.Item(35,9).Value = 'easy'
$ws.Cells.Item(36,9).Value = 'easy'
$ws.Cells.Item(37,9).Value = 'difficult'
$ws.Cells.Item(38,9).Value = 'difficult'
$ws.Cells.Item(39,9).Value = 'easy'
$ws.Cells.Item(40,9).Value = 'easy'
$ws.Cells.Item(41,9).Value = 'difficult'
$ws.Cells.Item(42,9).Value = 'easy'
$ws.Cells.Item(43,9).Value = 'easy'
$ws.Cells.Item(44,9).Value = 'easy'
$ws.Cells.Item(45,9).Value = 'difficult'
$ws.Cells.Item(46,9).Value = 'difficult'
$ws.Cells.Item(47,9).Value = 'easy'
$ws.Cells.Item(48,9).Value = 'easy'
$ws.Cells.Item(49,9).Value = 'easy'
$ws.Cells.Item(50,9).Value = 'easy'
$ws.Cells.Item(51,9).Value = 'difficult'
$ws.Cells.Item(52,9).Value = 'difficult'
$ws.Cells.Item(53,9).Value = 'extreme'
$ws.Cells.Item(54,9).Value = 'extreme'
$ws.Cells.Item(55,9).Value = 'difficult'
$ws.Cells.Item(56,9).Value = 'extreme'

# Fix typo: awarder -> awarded
$ws.Cells.Item(45,1).Value = "Which team was awarded with the 'FIFA Fair Play Trophy' in 2010?"

# IS HISTORIC / IS STATISTICAL / IS MODERN flags for each question row
$ws.Cells.Item(2,6).Value = 'no'
$ws.Cells.Item(2,7).Value = 'no'
$ws.Cells.Item(2,8).Value = 'yes'
$ws.Cells.Item(3,6).Value = 'no'
$ws.Cells.Item(3,7).Value = 'no'
$ws.Cells.Item(3,8).Value = 'yes'
$ws.Cells.Item(4,6).Value = 'no'
$ws.Cells.Item(4,7).Value = 'yes'
$ws.Cells.Item(4,8).Value = 'yes'
$ws.Cells.Item(5,6).Value = 'yes'
$ws.Cells.Item(5,7).Value = 'yes'
$ws.Cells.Item(5,8).Value = 'no'
$ws.Cells.Item(6,6).Value = 'no'
$ws.Cells.Item(6,7).Value = 'yes'
$ws.Cells.Item(6,8).Value = 'no'
$ws.Cells.Item(7,6).Value = 'yes'
$ws.Cells.Item(7,7).Value = 'yes'
$ws.Cells.Item(7,8).Value = 'no'
$ws.Cells.Item(8,6).Value = 'no'
$ws.Cells.Item(8,7).Value = 'no'
$ws.Cells.Item(8,8).Value = 'yes'
$ws.Cells.Item(9,6).Value = 'yes'
$ws.Cells.Item(9,7).Value = 'no'
$ws.Cells.Item(9,8).Value = 'no'
$ws.Cells.Item(10,6).Value = 'yes'
$ws.Cells.Item(10,7).Value = 'no'
$ws.Cells.Item(10,8).Value = 'no'
$ws.Cells.Item(11,6).Value = 'yes'
$ws.Cells.Item(11,7).Value = 'yes'
$ws.Cells.Item(11,8).Value = 'no'
$ws.Cells.Item(12,6).Value = 'yes'
$ws.Cells.Item(12,7).Value = 'no'
$ws.Cells.Item(12,8).Value = 'no'
$ws.Cells.Item(13,6).Value = 'yes'
$ws.Cells.Item(13,7).Value = 'yes'
$ws.Cells.Item(13,8).Value = 'no'
$ws.Cells.Item(14,6).Value = 'no'
$ws.Cells.Item(14,7).Value = 'yes'
$ws.Cells.Item(14,8).Value = 'no'
$ws.Cells.Item(15,6).Value = 'no'
$ws.Cells.Item(15,7).Value = 'yes'
$ws.Cells.Item(15,8).Value = 'no'
$ws.Cells.Item(16,6).Value = 'no'
$ws.Cells.Item(16,7).Value = 'yes'
$ws.Cells.Item(16,8).Value = 'no'
$ws.Cells.Item(17,6).Value = 'no'
$ws.Cells.Item(17,7).Value = 'no'
$ws.Cells.Item(17,8).Value = 'yes'
$ws.Cells.Item(18,6).Value = 'no'
$ws.Cells.Item(18,7).Value = 'yes'
$ws.Cells.Item(18,8).Value = 'no'
$ws.Cells.Item(19,6).Value = 'no'
$ws.Cells.Item(19,7).Value = 'yes'
$ws.Cells.Item(19,8).Value = 'no'
$ws.Cells.Item(20,6).Value = 'no'
$ws.Cells.Item(20,7).Value = 'yes'
$ws.Cells.Item(20,8).Value = 'no'
$ws.Cells.Item(21,6).Value = 'no'
$ws.Cells.Item(21,7).Value = 'yes'
$ws.Cells.Item(21,8).Value = 'no'
$ws.Cells.Item(22,6).Value = 'yes'
$ws.Cells.Item(22,7).Value = 'yes'
$ws.Cells.Item(22,8).Value = 'no'
$ws.Cells.Item(23,6).Value = 'no'
$ws.Cells.Item(23,7).Value = 'no'
$ws.Cells.Item(23,8).Value = 'yes'
$ws.Cells.Item(24,6).Value = 'no'
$ws.Cells.Item(24,7).Value = 'no'
$ws.Cells.Item(24,8).Value = 'yes'
$ws.Cells.Item(25,6).Value = 'yes'
$ws.Cells.Item(25,7).Value = 'no'
$ws.Cells.Item(25,8).Value = 'no'
$ws.Cells.Item(26,6).Value = 'yes'
$ws.Cells.Item(26,7).Value = 'no'
$ws.Cells.Item(26,8).Value = 'no'
$ws.Cells.Item(27,6).Value = 'yes'
$ws.Cells.Item(27,7).Value = 'no'
$ws.Cells.Item(27,8).Value = 'no'
$ws.Cells.Item(28,6).Value = 'no'
$ws.Cells.Item(28,7).Value = 'no'
$ws.Cells.Item(28,8).Value = 'yes'
$ws.Cells.Item(29,6).Value = 'no'
$ws.Cells.Item(29,7).Value = 'yes'
$ws.Cells.Item(29,8).Value = 'yes'
$ws.Cells.Item(30,6).Value = 'yes'
$ws.Cells.Item(30,7).Value = 'yes'
$ws.Cells.Item(30,8).Value = 'no'
$ws.Cells.Item(31,6).Value = 'no'
$ws.Cells.Item(31,7).Value = 'no'
$ws.Cells.Item(31,8).Value = 'yes'
$ws.Cells.Item(32,6).Value = 'no'
$ws.Cells.Item(32,7).Value = 'no'
$ws.Cells.Item(32,8).Value = 'yes'
$ws.Cells.Item(33,6).Value = 'no'
$ws.Cells.Item(33,7).Value = 'no'
$ws.Cells.Item(33,8).Value = 'yes'
$ws.Cells.Item(34,6).Value = 'no'
$ws.Cells.Item(34,7).Value = 'no'
$ws.Cells.Item(34,8).Value = 'yes'
$ws.Cells.Item(35,6).Value = 'no'
$ws.Cells.Item(35,7).Value = 'no'
$ws.Cells.Item(35,8).Value = 'yes'
$ws.Cells.Item(36,6).Value = 'no'
$ws.Cells.Item(36,7).Value = 'no'
$ws.Cells.Item(36,8).Value = 'yes'
$ws.Cells.Item(37,6).Value = 'no'
$ws.Cells.Item(37,7).Value = 'no'
$ws.Cells.Item(37,8).Value = 'yes'
$ws.Cells.Item(38,6).Value = 'no'
$ws.Cells.Item(38,7).Value = 'yes'
$ws.Cells.Item(38,8).Value = 'no'
$ws.Cells.Item(39,6).Value = 'no'
$ws.Cells.Item(39,7).Value = 'no'
$ws.Cells.Item(39,8).Value = 'yes'
$ws.Cells.Item(40,6).Value = 'no'
$ws.Cells.Item(40,7).Value = 'no'
$ws.Cells.Item(40,8).Value = 'yes'
$ws.Cells.Item(41,6).Value = 'no'
$ws.Cells.Item(41,7).Value = 'no'
$ws.Cells.Item(41,8).Value = 'yes'
$ws.Cells.Item(42,6).Value = 'no'
$ws.Cells.Item(42,7).Value = 'yes'
$ws.Cells.Item(42,8).Value = 'yes'
$ws.Cells.Item(43,6).Value = 'no'
$ws.Cells.Item(43,7).Value = 'yes'
$ws.Cells.Item(43,8).Value = 'yes'
$ws.Cells.Item(44,6).Value = 'no'
$ws.Cells.Item(44,7).Value = 'yes'
$ws.Cells.Item(44,8).Value = 'yes'
$ws.Cells.Item(45,6).Value = 'no'
$ws.Cells.Item(45,7).Value = 'yes'
$ws.Cells.Item(45,8).Value = 'yes'
$ws.Cells.Item(46,6).Value = 'no'
$ws.Cells.Item(46,7).Value = 'no'
$ws.Cells.Item(46,8).Value = 'yes'
$ws.Cells.Item(47,6).Value = 'no'
$ws.Cells.Item(47,7).Value = 'no'
$ws.Cells.Item(47,8).Value = 'yes'
$ws.Cells.Item(48,6).Value = 'no'
$ws.Cells.Item(48,7).Value = 'no'
$ws.Cells.Item(48,8).Value = 'yes'
$ws.Cells.Item(49,6).Value = 'no'
$ws.Cells.Item(49,7).Value = 'no'
$ws.Cells.Item(49,8).Value = 'yes'
$ws.Cells.Item(50,6).Value = 'no'
$ws.Cells.Item(50,7).Value = 'no'
$ws.Cells.Item(50,8).Value = 'yes'
$ws.Cells.Item(51,6).Value = 'no'
$ws.Cells.Item(51,7).Value = 'no'
$ws.Cells.Item(51,8).Value = 'yes'
$ws.Cells.Item(52,6).Value = 'yes'
$ws.Cells.Item(52,7).Value = 'no'
$ws.Cells.Item(52,8).Value = 'no'
$ws.Cells.Item(53,6).Value = 'yes'
$ws.Cells.Item(53,7).Value = 'no'
$ws.Cells.Item(53,8).Value = 'no'
$ws.Cells.Item(54,6).Value = 'yes'
$ws.Cells.Item(54,7).Value = 'no'
$ws.Cells.Item(54,8).Value = 'no'
$ws.Cells.Item(55,6).Value = 'yes'
$ws.Cells.Item(55,7).Value = 'no'
$ws.Cells.Item(55,8).Value = 'no'
$ws.Cells.Item(56,6).Value = 'yes'
$ws.Cells.Item(56,7).Value = 'no'
$ws.Cells.Item(56,8).Value = 'no'

# Column widths
$ws.Columns.Item(1).ColumnWidth = 115.140625
$ws.Columns.Item(9).ColumnWidth = 18.7109375

# Selection / view state
$ws.Range("I1:I68").Select()
